# Update stage 1 targetCurve: split the single "deviationPercent" column (I)
# into "deviationPercentX" (I) and a new "deviationPercentY" (J) column,
# carrying over the same numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Headers -----------------------------------------------------------
# Rename the existing header and add the new one.
$ws.Range("I1").Value = "deviationPercentX"
$ws.Range("J1").Value = "deviationPercentY"

# Give the new header the same formatting as the other numeric-column
# headers (e.g. B1/C1): the #,##0 number format (it already has the
# bordered, centered + wrapped header look, so just adding the number
# format makes it match that style exactly).
$ws.Range("I1").NumberFormat = "#,##0"

# --- Data row ------------------------------------------------------------
# Mirror the existing deviationPercent value into the new column, with the
# same numeric formatting I2 already has (bordered, centered, #,##0).
$ws.Range("J2").Value = $ws.Range("I2").Value()
$ws.Range("J2").NumberFormat = $ws.Range("I2").NumberFormat()
$ws.Range("J2").HorizontalAlignment = $ws.Range("I2").HorizontalAlignment()

# --- Column-wide formatting ----------------------------------------------
# The I column (rows 3:20, the unrelated lower block) keeps a plain,
# borderless style but now also carries the #,##0 number format so it
# matches the column's new numeric semantics.
$ws.Range("I3:I20").NumberFormat = "#,##0"
$ws.Range("I3:I20").HorizontalAlignment = 1

$ws.Columns.Item(9).NumberFormat = "#,##0"
$ws.Columns.Item(9).HorizontalAlignment = 1
